$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23; this shifts existing rows 23-34 down to 24-35
# (including their styles/formatting).
$ws.Range("A23").EntireRow.Insert()

# Populate the newly inserted row 23 with the new weekly price record.
$ws.Cells.Item(23, 1).Value  = 11
$ws.Cells.Item(23, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value  = "Bíobío"
$ws.Cells.Item(23, 4).Value  = 44572
$ws.Cells.Item(23, 5).Value  = 8
$ws.Cells.Item(23, 6).Value  = 100112031
$ws.Cells.Item(23, 7).Value  = "Poroto verde"
$ws.Cells.Item(23, 8).Value  = "Magnum"
$ws.Cells.Item(23, 9).Value  = "Primera"
$ws.Cells.Item(23, 10).Value = 250
$ws.Cells.Item(23, 11).Value = 26000
$ws.Cells.Item(23, 12).Value = 27000
$ws.Cells.Item(23, 13).Value = 26480
$ws.Cells.Item(23, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Región Metropolitana"
$ws.Cells.Item(23, 16).Value = 1059
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
